# BIS-769: Fixed xls test files
# Adds two new trailing columns ("Pattern" / "Pattern Type") to each of the
# three property-table headers (rows 4, 11, 19) in the sample-type export
# sheet, right after the existing "Unique" column (L), and moves the active
# selection to the newly added M19:N19 range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headerRows = @(4, 11, 19)

foreach ($row in $headerRows) {
    $srcCell = $ws.Range("L" + $row)
    $dstRange = $ws.Range("M" + $row + ":N" + $row)

    # Copy the "Unique" header cell formatting (bold header style) onto the
    # two new cells so they visually match the rest of the header row.
    $srcCell.Copy()
    $dstRange.PasteSpecial(-4122)

    $ws.Range("M" + $row).Value = "Pattern"
    $ws.Range("N" + $row).Value = "Pattern Type"
}

$excel.CutCopyMode = 0

# Match the author's final selection state recorded in the workbook.
$ws.Range("M19:N19").Select() | Out-Null
